# Apply updated probability-matrix values to Sheet1 (team_specific_matrix/Charlotte_B.xlsx)
# per "changes to team matrices from games pulled march 7"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Range("B2").Value = 0.2252559726962457
    $ws.Range("C2").Value = 0.4948805460750853
    $ws.Range("J2").Value = 0.01023890784982935
    $ws.Range("P2").Value = 0.1433447098976109
    $ws.Range("S2").Value = 0.1262798634812287
    $ws.Range("B3").Value = 0.01360544217687075
    $ws.Range("C3").Value = 0.02040816326530612
    $ws.Range("J3").Value = 0.02040816326530612
    $ws.Range("P3").Value = 0.8231292517006803
    $ws.Range("S3").Value = 0.1224489795918367
    $ws.Range("J4").Value = 0.08333333333333333
    $ws.Range("P4").Value = 0.6666666666666666
    $ws.Range("S4").Value = 0.25
    $ws.Range("B6").Value = 0.06334841628959276
    $ws.Range("D6").Value = 0.02714932126696833
    $ws.Range("F6").Value = 0.05429864253393665
    $ws.Range("J6").Value = 0.2398190045248869
    $ws.Range("O6").Value = 0.01809954751131222
    $ws.Range("Q6").Value = 0.1809954751131222
    $ws.Range("R6").Value = 0.07692307692307693
    $ws.Range("S6").Value = 0.3393665158371041
    $ws.Range("B7").Value = 0.102803738317757
    $ws.Range("D7").Value = 0.02336448598130841
    $ws.Range("F7").Value = 0.05607476635514019
    $ws.Range("J7").Value = 0.1682242990654206
    $ws.Range("O7").Value = 0.01869158878504673
    $ws.Range("Q7").Value = 0.191588785046729
    $ws.Range("R7").Value = 0.05607476635514019
    $ws.Range("S7").Value = 0.3831775700934579
    $ws.Range("B8").Value = 0.09401709401709402
    $ws.Range("D8").Value = 0.02777777777777778
    $ws.Range("E8").Value = 0.002136752136752137
    $ws.Range("F8").Value = 0.04700854700854701
    $ws.Range("J8").Value = 0.09401709401709402
    $ws.Range("O8").Value = 0.01923076923076923
    $ws.Range("Q8").Value = 0.1709401709401709
    $ws.Range("R8").Value = 0.1025641025641026
    $ws.Range("S8").Value = 0.4423076923076923
    $ws.Range("B9").Value = 0.07142857142857142
    $ws.Range("D9").Value = 0.02197802197802198
    $ws.Range("F9").Value = 0.03296703296703297
    $ws.Range("J9").Value = 0.1373626373626374
    $ws.Range("O9").Value = 0.02197802197802198
    $ws.Range("Q9").Value = 0.1813186813186813
    $ws.Range("R9").Value = 0.0989010989010989
    $ws.Range("S9").Value = 0.4340659340659341
    $ws.Range("B10").Value = 0.1027888446215139
    $ws.Range("D10").Value = 0.01752988047808765
    $ws.Range("F10").Value = 0.0796812749003984
    $ws.Range("J10").Value = 0.1266932270916335
    $ws.Range("O10").Value = 0.01195219123505976
    $ws.Range("Q10").Value = 0.2254980079681275
    $ws.Range("R10").Value = 0.07888446215139443
    $ws.Range("S10").Value = 0.3569721115537849
    $ws.Range("G11").Value = 0.1378205128205128
    $ws.Range("J11").Value = 0.0673076923076923
    $ws.Range("K11").Value = 0.1602564102564103
    $ws.Range("L11").Value = 0.6185897435897436
    $ws.Range("S11").Value = 0.01602564102564102
    $ws.Range("G12").Value = 0.7286432160804021
    $ws.Range("J12").Value = 0.2010050251256282
    $ws.Range("K12").Value = 0.01005025125628141
    $ws.Range("L12").Value = 0.03015075376884422
    $ws.Range("S12").Value = 0.03015075376884422
    $ws.Range("G13").Value = 0.673469387755102
    $ws.Range("J13").Value = 0.2448979591836735
    $ws.Range("S13").Value = 0.08163265306122448
    $ws.Range("F15").Value = 0.03153153153153153
    $ws.Range("H15").Value = 0.1576576576576577
    $ws.Range("I15").Value = 0.06756756756756757
    $ws.Range("J15").Value = 0.3603603603603603
    $ws.Range("K15").Value = 0.08558558558558559
    $ws.Range("M15").Value = 0.02702702702702703
    $ws.Range("O15").Value = 0.05855855855855856
    $ws.Range("S15").Value = 0.2117117117117117
    $ws.Range("F16").Value = 0.01595744680851064
    $ws.Range("H16").Value = 0.2340425531914894
    $ws.Range("I16").Value = 0.101063829787234
    $ws.Range("J16").Value = 0.3563829787234042
    $ws.Range("K16").Value = 0.101063829787234
    $ws.Range("M16").Value = 0.02127659574468085
    $ws.Range("O16").Value = 0.04787234042553191
    $ws.Range("S16").Value = 0.1223404255319149
    $ws.Range("F17").Value = 0.03151260504201681
    $ws.Range("H17").Value = 0.1827731092436975
    $ws.Range("I17").Value = 0.07142857142857142
    $ws.Range("J17").Value = 0.4012605042016807
    $ws.Range("K17").Value = 0.1176470588235294
    $ws.Range("M17").Value = 0.01470588235294118
    $ws.Range("O17").Value = 0.06722689075630252
    $ws.Range("S17").Value = 0.1134453781512605
    $ws.Range("F18").Value = 0.01570680628272251
    $ws.Range("H18").Value = 0.193717277486911
    $ws.Range("I18").Value = 0.06806282722513089
    $ws.Range("J18").Value = 0.4293193717277487
    $ws.Range("K18").Value = 0.08900523560209424
    $ws.Range("M18").Value = 0.01570680628272251
    $ws.Range("O18").Value = 0.1099476439790576
    $ws.Range("S18").Value = 0.07853403141361257
    $ws.Range("F19").Value = 0.01544715447154472
    $ws.Range("H19").Value = 0.216260162601626
    $ws.Range("I19").Value = 0.08211382113821138
    $ws.Range("J19").Value = 0.3691056910569105
    $ws.Range("K19").Value = 0.1203252032520325
    $ws.Range("M19").Value = 0.02520325203252033
    $ws.Range("N19").Value = 0.0008130081300813008
    $ws.Range("O19").Value = 0.06341463414634146
    $ws.Range("S19").Value = 0.1073170731707317
